# v1.1 - Added Processor selection support
$wb = $excel.ActiveWorkbook

# The second sheet ("getMobilePriceDetails") is the active / target sheet.
$ws = $wb.Worksheets.Item("getMobilePriceDetails")

# Add a new "Processor" column (F) with header + two values, mirroring the
# existing RAM-style column layout (header in row 1, values in rows 2-3).
$ws.Range("F1").Value = "Processor"
$ws.Range("F2").Value = "Snapdragon"
$ws.Range("F3").Value = "Exynos"

# Make the sheet active and select the newly added cell, matching the
# updated selection/activeCell in the workbook.
$ws.Activate()
$ws.Range("F3").Select()
